# issue #5: stock data output to json file
#
# Adds a "property_category" column (value "stock") to the 股票 (Stock)
# sheet, between the existing "total" and "date" columns, so the
# exported JSON records which property category each row belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 股票 (Stock) sheet

# Insert a new column at H, shifting date / legislator_name /
# legislator_id one column to the right (H->I, I->J, J->K).
$ws.Columns.Item(8).Insert()

# Header
$ws.Cells.Item(1, 8).Value = "property_category"

# Data rows
$ws.Cells.Item(2, 8).Value = "stock"
$ws.Cells.Item(3, 8).Value = "stock"
